$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading apostrophe forces text entry (avoids Excel auto-converting
# numeric-looking strings like '1.008' or '316.40' into real numbers),
# and resetting the Style back to Normal keeps the cell formatting
# identical to the original (no quote-prefix / text number format left behind).

$ws.Range("D2").Value = "'28.069.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "'1.908.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("D5").Value = "'316.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "'0.4831"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").Value = "'0.3820"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "'0.07364"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'0.9335"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").Value = "'20.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'0.07844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'1.887.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'5.499"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'6.607"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'91.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'0.000008823"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "'28.086.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'14.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "'5.152"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'2.150.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "'10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'156.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "'1.924"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "'18.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'2.098"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").Value = "'116.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'4.959"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'0.08902"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'3.362"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D34").Value = "'0.7666"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").Value = "'4.674"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "'2.603"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'0.02041"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "'1.096"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'0.05296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'0.5476"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "'2.987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").Value = "'7.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").Value = "'0.1522"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'8.445"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'0.4828"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "'106.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "'1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'1.654"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'68.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "'0.06095"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
